$wb = $excel.ActiveWorkbook

# Target B (Cutoff) and C (Reaction_number) values for rows 2-16 of each sheet.
# Column A (index 0..14) is unchanged; rows 17-20 are removed afterwards.

$nbrValues = @(
    @(5, 815),
    @(6, 815),
    @(7, 808),
    @(8, 774),
    @(9, 776),
    @(10, 778),
    @(11, 776),
    @(12, 767),
    @(13, 765),
    @(14, 759),
    @(15, 753),
    @(16, 753),
    @(17, 747),
    @(18, 750),
    @(19, 0)
)

$barValues = @(
    @(5, 723),
    @(6, 717),
    @(7, 722),
    @(8, 710),
    @(9, 703),
    @(10, 700),
    @(11, 704),
    @(12, 702),
    @(13, 702),
    @(14, 702),
    @(15, 702),
    @(16, 698),
    @(17, 697),
    @(18, 689),
    @(19, 0)
)

$wsNBR = $wb.Worksheets.Item("NBR")
$wsBAR = $wb.Worksheets.Item("BAR")

for ($i = 0; $i -lt $nbrValues.Count; $i++) {
    $row = $i + 2
    $pair = $nbrValues[$i]
    $wsNBR.Cells.Item($row, 2).Value = $pair[0]
    $wsNBR.Cells.Item($row, 3).Value = $pair[1]
}

for ($i = 0; $i -lt $barValues.Count; $i++) {
    $row = $i + 2
    $pair = $barValues[$i]
    $wsBAR.Cells.Item($row, 2).Value = $pair[0]
    $wsBAR.Cells.Item($row, 3).Value = $pair[1]
}

# Remove the now-obsolete trailing rows (17-20), which shrinks the sheet
# dimension from A1:C20 down to A1:C16.
$wsNBR.Range("A17:A20").EntireRow.Delete()
$wsBAR.Range("A17:A20").EntireRow.Delete()
